$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, [string]$Value)
    $origStyle = $Range.Style
    $Range.NumberFormat = "@"
    $Range.Value = $Value
    $Range.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "30.039.45"
Set-TextValue $ws.Range("E2") "  +0.04%  "

Set-TextValue $ws.Range("D3") "1.883.81"
Set-TextValue $ws.Range("E3") "  +0.45%  "

Set-TextValue $ws.Range("D4") "0.9985"
Set-TextValue $ws.Range("E4") "  -0.24%  "

Set-TextValue $ws.Range("D5") "243.95"
Set-TextValue $ws.Range("E5") "  -1.96%  "

Set-TextValue $ws.Range("D6") "0.9983"
Set-TextValue $ws.Range("E6") "  -0.28%  "

Set-TextValue $ws.Range("E7") "  -0.04%  "

Set-TextValue $ws.Range("D8") "44.34"
Set-TextValue $ws.Range("E8") "  -2.90%  "

Set-TextValue $ws.Range("D9") "0.2908"
Set-TextValue $ws.Range("E9") "  +2.29%  "

Set-TextValue $ws.Range("D10") "0.06620"
Set-TextValue $ws.Range("E10") "  +1.00%  "

Set-TextValue $ws.Range("D11") "1.881.45"
Set-TextValue $ws.Range("E11") "  +0.58%  "

Set-TextValue $ws.Range("D12") "16.80"
Set-TextValue $ws.Range("E12") "  -1.58%  "

Set-TextValue $ws.Range("D13") "0.07191"
Set-TextValue $ws.Range("E13") "  +0.07%  "

Set-TextValue $ws.Range("D14") "0.6644"
Set-TextValue $ws.Range("E14") "  +0.57%  "

Set-TextValue $ws.Range("D15") "85.70"
Set-TextValue $ws.Range("E15") "  +0.75%  "

Set-TextValue $ws.Range("D16") "4.851"
Set-TextValue $ws.Range("E16") "  +1.09%  "

Set-TextValue $ws.Range("D17") "30.040.75"
Set-TextValue $ws.Range("E17") "  +0.08%  "

Set-TextValue $ws.Range("D18") "0.000007774"
Set-TextValue $ws.Range("E18") "  +3.57%  "

Set-TextValue $ws.Range("D19") "0.9983"
Set-TextValue $ws.Range("E19") "  -0.05%  "

Set-TextValue $ws.Range("D20") "12.76"
Set-TextValue $ws.Range("E20") "  -0.82%  "

Set-TextValue $ws.Range("D21") "2.121.08"

Set-TextValue $ws.Range("D22") "0.9979"
Set-TextValue $ws.Range("E22") "  -0.26%  "

Set-TextValue $ws.Range("D23") "4.762"
Set-TextValue $ws.Range("E23") "  +0.36%  "

Set-TextValue $ws.Range("D24") "5.599"
Set-TextValue $ws.Range("E24") "  +1.77%  "

Set-TextValue $ws.Range("D25") "9.151"
Set-TextValue $ws.Range("E25") "  +1.45%  "

Set-TextValue $ws.Range("D26") "150.35"
Set-TextValue $ws.Range("E26") "  +4.12%  "

Set-TextValue $ws.Range("D27") "135.87"
Set-TextValue $ws.Range("E27") "  +0.91%  "

Set-TextValue $ws.Range("D28") "16.77"
Set-TextValue $ws.Range("E28") "  +0.32%  "

Set-TextValue $ws.Range("E29") "  -2.46%  "

Set-TextValue $ws.Range("D30") "1.378"
Set-TextValue $ws.Range("E30") "  -0.52%  "

Set-TextValue $ws.Range("D31") "4.169"
Set-TextValue $ws.Range("E31") "  -0.87%  "

Set-TextValue $ws.Range("D32") "0.08673"
Set-TextValue $ws.Range("E32") "  +0.90%  "

Set-TextValue $ws.Range("D33") "3.947"
Set-TextValue $ws.Range("E33") "  +1.69%  "

Set-TextValue $ws.Range("D34") "0.04995"
Set-TextValue $ws.Range("E34") "  -1.40%  "

Set-TextValue $ws.Range("D35") "1.105"
Set-TextValue $ws.Range("E35") "  -2.91%  "

Set-TextValue $ws.Range("D36") "0.7030"
Set-TextValue $ws.Range("E36") "  +2.76%  "

Set-TextValue $ws.Range("D37") "2.656"
Set-TextValue $ws.Range("E37") "  -1.77%  "

$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D38") "2.699"
Set-TextValue $ws.Range("E38") "  -1.54%  "

$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D39") "2.193"
Set-TextValue $ws.Range("E39") "  -5.55%  "

Set-TextValue $ws.Range("D40") "0.9347"
Set-TextValue $ws.Range("E40") "  -2.84%  "

Set-TextValue $ws.Range("E41") "  +1.15%  "

Set-TextValue $ws.Range("D42") "5.953"
Set-TextValue $ws.Range("E42") "  -1.95%  "

Set-TextValue $ws.Range("D43") "0.9990"
Set-TextValue $ws.Range("E43") "  -0.20%  "

Set-TextValue $ws.Range("D44") "0.4185"
Set-TextValue $ws.Range("E44") "  -0.07%  "

Set-TextValue $ws.Range("D45") "101.32"
Set-TextValue $ws.Range("E45") "  -1.72%  "

Set-TextValue $ws.Range("D46") "7.486"
Set-TextValue $ws.Range("E46") "  -0.16%  "

Set-TextValue $ws.Range("D47") "0.1261"
Set-TextValue $ws.Range("E47") "  +0.67%  "

Set-TextValue $ws.Range("D48") "0.05716"
Set-TextValue $ws.Range("E48") "  +1.62%  "

Set-TextValue $ws.Range("D49") "32.36"
Set-TextValue $ws.Range("E49") "  -0.21%  "

Set-TextValue $ws.Range("D50") "8.252"
Set-TextValue $ws.Range("E50") "  +0.57%  "

$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D51") "1.341"
Set-TextValue $ws.Range("E51") "  +0.39%  "
